$d = $word.ActiveDocument

$replacements = @(
    @("261×6=1566", "160×7=1120"),
    @("101×6=606", "110×4=440"),
    @("677×4=2708", "362×2=724"),
    @("661×6=3966", "535×9=4815"),
    @("276×7=1932", "372×5=1860"),
    @("220×6=1320", "854×2=1708"),
    @("658×2=1316", "266×4=1064"),
    @("607×5=3035", "247×6=1482"),
    @("953×4=3812", "135×5=675"),
    @("967×5=4835", "913×9=8217"),
    @("268×5=1340", "944×3=2832"),
    @("743×5=3715", "272×3=816"),
    @("864×4=3456", "355×5=1775"),
    @("642×8=5136", "475×2=950"),
    @("899×6=5394", "432×9=3888"),
    @("320×3=960", "532×5=2660"),
    @("653×5=3265", "429×9=3861"),
    @("137×2=274", "110×8=880"),
    @("633×7=4431", "947×5=4735"),
    @("947×8=7576", "270×2=540"),
    @("959×9=8631", "708×8=5664"),
    @("767×5=3835", "350×9=3150"),
    @("775×2=1550", "203×8=1624"),
    @("419×7=2933", "298×8=2384"),
    @("372×8=2976", "218×6=1308")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
